# Remove the trailing "Ver no Jupiter..." / copyright footer block that was
# added by the site build, along with the blank paragraph that separated it
# from the bibliography entry. The blank paragraph that follows the footer
# (right before the page-break paragraph) is left untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph using Find (content-driven, not
# index-driven, so it's robust to any paragraph renumbering).
$finder = $d.Content
$found = $finder.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Ver no Jupiter' paragraph"
}

$targetStart = $finder.Start

# Map that character offset back to a paragraph index.
$jupiterIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $targetStart) {
        $jupiterIndex = $i
        break
    }
}

if ($jupiterIndex -eq -1) {
    throw "Could not resolve paragraph index for the 'Ver no Jupiter' match"
}

# Delete the empty paragraph right before it, the "Ver no Jupiter..."
# paragraph itself, and the "(c) 2020 ..." paragraph right after it.
$startPara = $d.Paragraphs.Item($jupiterIndex - 1)
$endPara = $d.Paragraphs.Item($jupiterIndex + 1)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
